# Apply the commit's content edits to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "product_id" -> "id"
$ws.Range("A1").Value = "id"

# Row 2 (was product_id 3): renumber id, translate title, bump price
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "襯衫"
$ws.Range("E2").Value = 30

# Row 3 (was product_id 4): renumber id, translate title, extend color_ids
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = "另一個裙子"
$ws.Range("L3").Value = "2,13"

# Move the active selection to D15, matching the saved view state.
$ws.Range("D15").Select()
